# Automated daily data refresh: appends the new trading day
# (05-sep for the "Prix Spot" hourly matrix, 2025-09-03 for the
# "Gaz" and "CO2" daily series) to the three worksheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Prix Spot": new column CF ("05-sep") with hourly prices
# ---------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the header style (bold, centered, bordered) from the previous
# day's header cell (CE1) onto the new header cell (CF1), then set
# its value.
$wsSpot.Range("CE1").Copy($wsSpot.Range("CF1"))
$wsSpot.Range("CF1").Value = "05-sep"

$spotValues = @{
    2  = 70.86
    3  = 58.98
    4  = 61.49
    5  = 50.12
    6  = 52.65
    7  = 52.53
    8  = 76.79
    9  = 89.99
    10 = 88.22
    11 = 79.1
    12 = 35
    13 = 10
    14 = 15
    15 = 18.43
    16 = 6.62
    17 = 17.07
    18 = 18.63
    19 = 30.4
    20 = 45.23
    21 = 89.25
    22 = 103.97
    23 = 107.1
    24 = 95
    25 = 88.38
}

foreach ($row in $spotValues.Keys) {
    $wsSpot.Cells.Item($row, 84).Value = $spotValues[$row]
}

# ---------------------------------------------------------------
# Sheet "Gaz": append 2025-09-03 row
# ---------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force text formatting before assigning the date-like string so the
# engine keeps it as literal text ("2025-09-03") instead of coercing
# it into a date serial, then drop the temporary number format so the
# new cell stays unstyled, matching the rest of column A.
$wsGaz.Range("A81").NumberFormat = "@"
$wsGaz.Range("A81").Value = "2025-09-03"
$wsGaz.Range("A81").ClearFormats()

$wsGaz.Range("B81").Value = 30.875

# ---------------------------------------------------------------
# Sheet "CO2": append 2025-09-03 row (settlement price not yet
# published, so the price cell is left blank)
# ---------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A81").NumberFormat = "@"
$wsCo2.Range("A81").Value = "2025-09-03"
$wsCo2.Range("A81").ClearFormats()

$wsCo2.Range("B81").NumberFormat = "@"
$wsCo2.Range("B81").Value = ""
$wsCo2.Range("B81").ClearFormats()
